$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(3)

# Add the new "Comments" header in column D
$ws.Range("D1").Value = "Comments"

# Add index numbers in column A for rows 3-5
$ws.Range("A3").Value = 1
$ws.Range("A4").Value = 2
$ws.Range("A5").Value = 3

# Move the comment strings from column B to column D
$ws.Range("D3").Value = $ws.Range("B3").Value2
$ws.Range("D4").Value = $ws.Range("B4").Value2
$ws.Range("D5").Value = $ws.Range("B5").Value2

$ws.Range("B3").Clear()
$ws.Range("B4").Clear()
$ws.Range("B5").Clear()

# Update the active selection to match the new layout
$ws.Activate() | Out-Null
$ws.Range("D3:D5").Select() | Out-Null
